$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for ticker "CDAY" (Ceridian) was removed. All rows below it, up
# through the row preceding "DE" (John Deere), shift up by one row. A new
# row is appended at the end of that block for ticker "DAY" (Dayforce),
# reusing the original Ceridian row's other attributes (sector, sub-sector,
# headquarters, date added, CIK, founded year).

$firstRow = 101
$lastRow = 143
$numCols = 8

# Columns whose text must be preserved verbatim (leading zeros, date-like
# strings, etc.) rather than being auto-coerced into numbers/dates by Excel.
# Prefixing with an apostrophe forces Excel to store the value as text.
function ForceText($val) {
    if ($val -eq $null) { return $val }
    return "'" + $val
}

# Capture the original rows (firstRow+1)..lastRow (columns A..H) before
# overwriting anything, since we shift them up by one row.
$buffer = @()
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value()
    }
    $buffer += ,$rowVals
}

# Write the captured rows shifted up by one (into firstRow..lastRow-1).
for ($i = 0; $i -lt $buffer.Length; $i++) {
    $targetRow = $firstRow + $i
    $rowVals = $buffer[$i]
    for ($c = 1; $c -le $numCols; $c++) {
        if ($c -ge 6) {
            $ws.Cells.Item($targetRow, $c).Value = ForceText $rowVals[$c - 1]
        } else {
            $ws.Cells.Item($targetRow, $c).Value = $rowVals[$c - 1]
        }
    }
}

# Write the new Dayforce row at the end of the block, keeping the same
# sector/sub-sector/headquarters/date-added/CIK/founded values that the
# Ceridian row originally had.
$ws.Cells.Item($lastRow, 1).Value = "DAY"
$ws.Cells.Item($lastRow, 2).Value = "Dayforce"
$ws.Cells.Item($lastRow, 3).Value = "Industrials"
$ws.Cells.Item($lastRow, 4).Value = "Human Resource & Employment Services"
$ws.Cells.Item($lastRow, 5).Value = "Minneapolis, Minnesota"
$ws.Cells.Item($lastRow, 6).Value = "'2021-09-20"
$ws.Cells.Item($lastRow, 7).Value = "'0001725057"
$ws.Cells.Item($lastRow, 8).Value = "'1992"

# Unrelated sub-sector rename: "Specialty Stores" -> "Other Specialty Retail"
# for Tractor Supply (row 448) and Ulta Beauty (row 459).
$ws.Cells.Item(448, 4).Value = "Other Specialty Retail"
$ws.Cells.Item(459, 4).Value = "Other Specialty Retail"
